# "add new function delete product"
# Deletes product id 1 (Iphone5) from the "Sheet" inventory sheet, which
# shifts the remaining product (Iphone8) up one row, and appends a log
# entry to the "Logs" sheet recording the deletion.

$wb = $excel.ActiveWorkbook

# --- 1. Delete the product row (Product ID 1 / Iphone5) from "Sheet" ---
$ws = $wb.Worksheets.Item("Sheet")
$ws.Rows.Item(2).Delete()

# --- 2. Append a "Deleted" entry to the "Logs" sheet ---
$logs = $wb.Worksheets.Item("Logs")
$logRow = $logs.UsedRange.Rows.Count + 1

# Leading "'" forces text (matches the existing plain-text Date /
# Product ID columns) instead of Excel auto-converting to a date/number.
$logs.Cells.Item($logRow, 1).Value = "'2025-04-20"
$logs.Cells.Item($logRow, 2).Value = "'1"
$logs.Cells.Item($logRow, 3).Value = "-"
$logs.Cells.Item($logRow, 4).Value = "-"
$logs.Cells.Item($logRow, 5).Value = "Deleted"
